$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 166.66667
$ws.Cells.Item(5, 9).Value = 166.66667
$ws.Cells.Item(5, 11).Value = 166.66667
$ws.Cells.Item(5, 13).Value = -51.66667000000001

$ws.Cells.Item(28, 8).Value = 3242.25
$ws.Cells.Item(28, 9).Value = 1046.25
$ws.Cells.Item(28, 10).Value = 4340.25
$ws.Cells.Item(28, 11).Value = 1046.25
$ws.Cells.Item(28, 12).Value = 4340.25
$ws.Cells.Item(28, 13).Value = -561.25
$ws.Cells.Item(28, 14).Value = -5310.25

$ws.Cells.Item(32, 8).Value = 500499.5
$ws.Cells.Item(32, 9).Value = 1000
$ws.Cells.Item(32, 10).Value = 999999
$ws.Cells.Item(32, 11).Value = 1000
$ws.Cells.Item(32, 12).Value = 999999
$ws.Cells.Item(32, 13).Value = -674
$ws.Cells.Item(32, 14).Value = -1000651

$ws.Cells.Item(33, 8).Value = 416.16666
$ws.Cells.Item(33, 9).Value = 329.4
$ws.Cells.Item(33, 10).Value = 850
$ws.Cells.Item(33, 11).Value = 329.4
$ws.Cells.Item(33, 12).Value = 850
$ws.Cells.Item(33, 13).Value = -100.4
$ws.Cells.Item(33, 14).Value = -1308

$ws.Cells.Item(40, 8).Value = 7570.15
$ws.Cells.Item(40, 9).Value = 3667.3333
$ws.Cells.Item(40, 10).Value = 10763.363
$ws.Cells.Item(40, 11).Value = 3667.3333
$ws.Cells.Item(40, 12).Value = 10763.363
$ws.Cells.Item(40, 13).Value = -3492.3333
$ws.Cells.Item(40, 14).Value = -11113.363

$ws.Cells.Item(64, 8).Value = 7999.8335
$ws.Cells.Item(64, 9).Value = 7999.8335
$ws.Cells.Item(64, 11).Value = 7999.8335
$ws.Cells.Item(64, 13).Value = -7751.8335

$ws.Cells.Item(67, 8).Value = 7999.8335
$ws.Cells.Item(67, 9).Value = 7999.8335
$ws.Cells.Item(67, 11).Value = 7999.8335
$ws.Cells.Item(67, 13).Value = -7141.8335

$ws.Cells.Item(100, 8).Value = 3940.7144
$ws.Cells.Item(100, 9).Value = 3747.6667
$ws.Cells.Item(100, 10).Value = 5099
$ws.Cells.Item(100, 11).Value = 3747.6667
$ws.Cells.Item(100, 12).Value = 5099
$ws.Cells.Item(100, 13).Value = -3206.6667
$ws.Cells.Item(100, 14).Value = -6181

$ws.Cells.Item(116, 8).Value = 8972.166999999999
$ws.Cells.Item(116, 9).Value = 7879.5
$ws.Cells.Item(116, 10).Value = 9518.5
$ws.Cells.Item(116, 11).Value = 7879.5
$ws.Cells.Item(116, 12).Value = 9518.5
$ws.Cells.Item(116, 13).Value = -4437.5
$ws.Cells.Item(116, 14).Value = -16402.5

$ws.Cells.Item(134, 8).Value = 52218.668
$ws.Cells.Item(134, 10).Value = 52218.668
$ws.Cells.Item(134, 12).Value = 52218.668
$ws.Cells.Item(134, 14).Value = -62358.668

$ws.Cells.Item(135, 8).Value = 1486.9375
$ws.Cells.Item(135, 9).Value = 1226.7693
$ws.Cells.Item(135, 10).Value = 2614.3333
$ws.Cells.Item(135, 11).Value = 11040.9237
$ws.Cells.Item(135, 12).Value = 23528.9997
$ws.Cells.Item(135, 13).Value = -8505.923699999999
$ws.Cells.Item(135, 14).Value = -28598.9997

$ws.Cells.Item(137, 8).Value = 727721.4
$ws.Cells.Item(137, 9).Value = 2354.625
$ws.Cells.Item(137, 11).Value = 7063.875
$ws.Cells.Item(137, 13).Value = -4513.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3726.7356
$ws.Cells.Item(32, 9).Value = 1343.7465
$ws.Cells.Item(32, 11).Value = 1343.7465
$ws.Cells.Item(32, 13).Value = -1056.7465

$ws.Cells.Item(45, 8).Value = 14002824
$ws.Cells.Item(45, 9).Value = 3356
$ws.Cells.Item(45, 10).Value = 25202398
$ws.Cells.Item(45, 11).Value = 3356
$ws.Cells.Item(45, 12).Value = 25202398
$ws.Cells.Item(45, 13).Value = -2979
$ws.Cells.Item(45, 14).Value = -25203152

$ws.Cells.Item(61, 8).Value = 33515.03
$ws.Cells.Item(61, 9).Value = 1109.8889
$ws.Cells.Item(61, 11).Value = 1109.8889
$ws.Cells.Item(61, 13).Value = -897.8888999999999

$ws.Cells.Item(74, 8).Value = 79855.30499999999
$ws.Cells.Item(74, 9).Value = 126771.25
$ws.Cells.Item(74, 10).Value = 4789.8
$ws.Cells.Item(74, 11).Value = 126771.25
$ws.Cells.Item(74, 12).Value = 4789.8
$ws.Cells.Item(74, 13).Value = -125897.25
$ws.Cells.Item(74, 14).Value = -6537.8

$ws.Cells.Item(77, 8).Value = 79855.30499999999
$ws.Cells.Item(77, 9).Value = 126771.25
$ws.Cells.Item(77, 10).Value = 4789.8
$ws.Cells.Item(77, 11).Value = 633856.25
$ws.Cells.Item(77, 12).Value = 23949
$ws.Cells.Item(77, 13).Value = -629488.25
$ws.Cells.Item(77, 14).Value = -32685

$ws.Cells.Item(88, 8).Value = 627.8
$ws.Cells.Item(88, 9).Value = 780
$ws.Cells.Item(88, 10).Value = 399.5
$ws.Cells.Item(88, 11).Value = 780
$ws.Cells.Item(88, 12).Value = 399.5
$ws.Cells.Item(88, 13).Value = -374
$ws.Cells.Item(88, 14).Value = -1211.5

$ws.Cells.Item(91, 8).Value = 627.8
$ws.Cells.Item(91, 9).Value = 780
$ws.Cells.Item(91, 10).Value = 399.5
$ws.Cells.Item(91, 11).Value = 780
$ws.Cells.Item(91, 12).Value = 399.5
$ws.Cells.Item(91, 13).Value = 624
$ws.Cells.Item(91, 14).Value = -3207.5

$ws.Cells.Item(107, 8).Value = 72966.8
$ws.Cells.Item(107, 10).Value = 72966.8
$ws.Cells.Item(107, 12).Value = 72966.8
$ws.Cells.Item(107, 14).Value = -80646.8

$ws.Cells.Item(132, 8).Value = 2370.5557
$ws.Cells.Item(132, 9).Value = 2081.8333
$ws.Cells.Item(132, 10).Value = 2948
$ws.Cells.Item(132, 11).Value = 6245.499899999999
$ws.Cells.Item(132, 12).Value = 8844
$ws.Cells.Item(132, 13).Value = -3715.499899999999
$ws.Cells.Item(132, 14).Value = -13904

$ws.Cells.Item(136, 8).Value = 33515.03
$ws.Cells.Item(136, 9).Value = 1109.8889
$ws.Cells.Item(136, 11).Value = 3329.6667
$ws.Cells.Item(136, 13).Value = -779.6666999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1213.375
$ws.Cells.Item(20, 9).Value = 1311.6
$ws.Cells.Item(20, 11).Value = 1311.6
$ws.Cells.Item(20, 13).Value = -1064.6

$ws.Cells.Item(22, 8).Value = 81166.16
$ws.Cells.Item(22, 9).Value = 104848.7
$ws.Cells.Item(22, 11).Value = 104848.7
$ws.Cells.Item(22, 13).Value = -104675.7

$ws.Cells.Item(94, 8).Value = 2631.0527
$ws.Cells.Item(94, 9).Value = 2617.0588
$ws.Cells.Item(94, 11).Value = 2617.0588
$ws.Cells.Item(94, 13).Value = -2166.0588

$ws.Cells.Item(134, 8).Value = 4406.9165
$ws.Cells.Item(134, 9).Value = 2893.6191
$ws.Cells.Item(134, 11).Value = 8680.8573
$ws.Cells.Item(134, 13).Value = -6145.8573

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(29, 8).Value = 1966.3334
$ws.Cells.Item(29, 10).Value = 1949.5
$ws.Cells.Item(29, 12).Value = 1949.5
$ws.Cells.Item(29, 14).Value = -2535.5

$ws.Cells.Item(32, 8).Value = 549.8333
$ws.Cells.Item(32, 9).Value = 562.5
$ws.Cells.Item(32, 10).Value = 524.5
$ws.Cells.Item(32, 11).Value = 562.5
$ws.Cells.Item(32, 12).Value = 524.5
$ws.Cells.Item(32, 13).Value = -246.5
$ws.Cells.Item(32, 14).Value = -1156.5

$ws.Cells.Item(35, 8).Value = 3130.9375
$ws.Cells.Item(35, 9).Value = 3599.5833
$ws.Cells.Item(35, 11).Value = 3599.5833
$ws.Cells.Item(35, 13).Value = -3305.5833

$ws.Cells.Item(59, 8).Value = 142474.25
$ws.Cells.Item(59, 10).Value = 89965.336
$ws.Cells.Item(59, 12).Value = 89965.336
$ws.Cells.Item(59, 14).Value = -92255.336

$ws.Cells.Item(107, 8).Value = 2004.125
$ws.Cells.Item(107, 10).Value = 2368.5
$ws.Cells.Item(107, 12).Value = 2368.5
$ws.Cells.Item(107, 14).Value = -6208.5

$ws.Cells.Item(132, 8).Value = 1625435.2
$ws.Cells.Item(132, 9).Value = 1625435.2
$ws.Cells.Item(132, 11).Value = 4876305.6
$ws.Cells.Item(132, 13).Value = -4873775.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 77180.62
$ws.Cells.Item(6, 9).Value = 91076.91
$ws.Cells.Item(6, 10).Value = 751
$ws.Cells.Item(6, 11).Value = 273230.73
$ws.Cells.Item(6, 12).Value = 2253
$ws.Cells.Item(6, 13).Value = -273117.73
$ws.Cells.Item(6, 14).Value = -2479

$ws.Cells.Item(56, 8).Value = 5956.4443
$ws.Cells.Item(56, 9).Value = 5956.4443
$ws.Cells.Item(56, 11).Value = 5956.4443
$ws.Cells.Item(56, 13).Value = -5426.4443

$ws.Cells.Item(98, 8).Value = 11001
$ws.Cells.Item(98, 9).Value = 14003
$ws.Cells.Item(98, 10).Value = 7999
$ws.Cells.Item(98, 11).Value = 42009
$ws.Cells.Item(98, 12).Value = 23997
$ws.Cells.Item(98, 13).Value = -40511
$ws.Cells.Item(98, 14).Value = -26993

$ws.Cells.Item(132, 8).Value = 4478.8887
$ws.Cells.Item(132, 10).Value = 6210.8335
$ws.Cells.Item(132, 12).Value = 55897.5015
$ws.Cells.Item(132, 14).Value = -60957.5015

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 996.6429000000001
$ws.Cells.Item(97, 9).Value = 458.15384
$ws.Cells.Item(97, 11).Value = 458.15384
$ws.Cells.Item(97, 13).Value = 37.84616

$ws.Cells.Item(122, 8).Value = 5623.3335
$ws.Cells.Item(122, 9).Value = 6705.8887
$ws.Cells.Item(122, 11).Value = 20117.6661
$ws.Cells.Item(122, 13).Value = -17667.6661

$ws.Cells.Item(132, 8).Value = 5288.875
$ws.Cells.Item(132, 9).Value = 4337.2
$ws.Cells.Item(132, 11).Value = 13011.6
$ws.Cells.Item(132, 13).Value = -10481.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 885.17645
$ws.Cells.Item(22, 10).Value = 854.8
$ws.Cells.Item(22, 12).Value = 854.8
$ws.Cells.Item(22, 14).Value = -1444.8

$ws.Cells.Item(27, 8).Value = 885.17645
$ws.Cells.Item(27, 10).Value = 854.8
$ws.Cells.Item(27, 12).Value = 854.8
$ws.Cells.Item(27, 14).Value = -1068.8

$ws.Cells.Item(40, 8).Value = 12967280
$ws.Cells.Item(40, 9).Value = 5725.25
$ws.Cells.Item(40, 11).Value = 5725.25
$ws.Cells.Item(40, 13).Value = -5589.25

$ws.Cells.Item(46, 8).Value = 2117.3333
$ws.Cells.Item(46, 9).Value = 970.1429000000001
$ws.Cells.Item(46, 11).Value = 970.1429000000001
$ws.Cells.Item(46, 13).Value = -782.1429000000001

$ws.Cells.Item(61, 8).Value = 2247.5
$ws.Cells.Item(61, 9).Value = 2247.5
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 2247.5
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -2045.5
$ws.Cells.Item(61, 14).ClearContents()

$ws.Cells.Item(113, 8).Value = 2247.5
$ws.Cells.Item(113, 9).Value = 2247.5
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 2247.5
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -77.5
$ws.Cells.Item(113, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1451.6207
$ws.Cells.Item(132, 9).Value = 1167.1923
$ws.Cells.Item(132, 11).Value = 3501.5769
$ws.Cells.Item(132, 13).Value = -971.5769
